# xleash: AMIDST impl p-eval filter.
# Recreate the authored edit against tests/recursive.xlsx:
#  - bump workbook window chrome (best effort; host may not persist this)
#  - add a new "eval sheet" worksheet (4th tab) with P-eval sample data
#  - insert a new row 9 on sheet 1 holding the P-eval trigger cell
#  - re-select/activate so the final UI state matches the authored file

$wb = $excel.ActiveWorkbook

# --- 1. Workbook window chrome -------------------------------------------
# (xWindow/yWindow/tabRatio move a bit in the authored diff; best effort.)
$win = $wb.Windows.Item(1)
$win.Left = 2265
$win.Top = 1335
$win.TabRatio = 0.358

$ws1 = $wb.Worksheets.Item(1)

# --- 2. sheet1: make room for the new "P-eval" trigger row ---------------
# Old row 10 (".. recurse" probe) through row 15 shift down by one; a brand
# new row 9 appears holding the eval-sheet trigger formula text.
$ws1.Rows.Item(9).Insert() | Out-Null
$ws1.Range("B9").Value = "P-eval"

# --- 3. Add the new "eval sheet" worksheet after the last tab ------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsEval = $wb.Worksheets.Add($null, $lastSheet)
$wsEval.Name = "eval sheet"

# Fill the sample data (order mirrors how the sheet was authored).
$wsEval.Range("C4").Value = "dict(a_dict=1)"
$wsEval.Range("C1").Value = "EVAL_COL"
$wsEval.Range("D1").Value = "NO_EVAL"

# Row 2: B2/D2 are quote-prefixed text (leading apostrophe forces text and
# stamps the quotePrefix="1" style, matching the authored cellXfs entry).
$wsEval.Range("D2").Value = "'a'+4"
$wsEval.Range("D3").Value = "bad boy"
$wsEval.Range("C2").Value = "a=1; a+5"
$wsEval.Range("C3").Value = "[1,2,3]"
$wsEval.Range("B4").Value = "bus"

# Remaining header/reused-string cells.
$wsEval.Range("B1").Value = "COL1"
$wsEval.Range("B2").Value = "'foo"
$wsEval.Range("B3").Value = "bar"

# Leftover selection state on the new sheet (matches authored file).
$wsEval.Range("D13").Select() | Out-Null

# --- 4. sheet1: write the P-eval trigger formula text last ---------------
# (Keeps shared-string allocation order matching the authored file.)
$ws1.Range("C9").Value = '#eval sheet!::{"opts": {"lax": true}, "func": "pipe", "args":[["df", {"index_col": null}]]}'

# --- 5. Restore sheet1 as the active tab/selection ------------------------
$ws1.Activate() | Out-Null
$ws1.Range("C9").Select() | Out-Null
